$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B, shifting the existing "Functional Location" and
# "Name" columns one place to the right (B->C, C->D).
$ws.Columns("B:B").Insert()

# New column B is "Object Type", populated with "PODE" for every data row.
$ws.Range("B1").Value = "Object Type"

# Column A keeps its position but becomes the "S4 Equipment Id" column with
# sequential ids (the original sample data repeated the same id on every
# row - fix that up too).
for ($r = 2; $r -le 11; $r++) {
  $ws.Cells.Item($r, 1).Value = 1000101000 + ($r - 1)
  $ws.Cells.Item($r, 2).Value = "PODE"
}

$ws.Range("A1").Value = "S4 Equipment Id"

# Give the newly inserted column a sensible width (matches column A).
$ws.Columns(2).ColumnWidth = 21.67

# Restore the active selection to A2 (frozen pane's top-left data cell).
$ws.Range("A2").Select() | Out-Null
